$d = $word.ActiveDocument

# --- Edit 1: first paragraph -------------------------------------------
# "This is a Microsoft word document." -> add two trailing spaces, then
# append a red-colored parenthetical note split across three runs.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute("This is a Microsoft word document.", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "This is a Microsoft word document.  ", 2)

$redColor = 192  # RGB(192,0,0) == C00000 as a Word OLE color value

$insPos = $p1.Range.End - 1
$ins = $d.Range($insPos, $insPos)
$ins.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$seg = $d.Range($insPos, $p1.Range.End - 1)
$seg.Font.Color = $redColor

$insPos2 = $p1.Range.End - 1
$ins2 = $d.Range($insPos2, $insPos2)
$ins2.InsertAfter("rsion for branch alternate")
$seg2 = $d.Range($insPos2, $p1.Range.End - 1)
$seg2.Font.Color = $redColor

$insPos3 = $p1.Range.End - 1
$ins3 = $d.Range($insPos3, $insPos3)
$ins3.InsertAfter(")")
$seg3 = $d.Range($insPos3, $p1.Range.End - 1)
$seg3.Font.Color = $redColor

# --- Edit 2: new shaded paragraph after the final paragraph -------------
# Re-assigning the trailing paragraph's own Range.Text (rather than using
# InsertParagraphAfter, which leaves a stray empty run behind, or Find's
# Execute/Replace, which rewrites the preceding run) appends a clean new
# paragraph mark without introducing any extra run into the new paragraph.
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$fullText = $pLast.Range.Text
$textNoMark = $fullText.Substring(0, $fullText.Length - 1)
$pStart = $pLast.Range.Start
$pEnd = $pLast.Range.End - 1
$rngText = $d.Range($pStart, $pEnd)
$rngText.Text = $textNoMark + "`r"

$newLastIndex = $d.Paragraphs.Count
$newP = $d.Paragraphs.Item($newLastIndex)
$newP.Range.Style = "Normal"
$newP.Format.Shading.Texture = 0
$newP.Format.Shading.ForegroundPatternColor = -16777216
$newP.Format.Shading.BackgroundPatternColor = 16382457
